$d = $word.ActiveDocument

# 1. Extend the "What income bracket..." question with the extra follow-up text.
$d.Content.Find.Execute(
    "What income bracket are residential users in?",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "What income bracket are residential users in? Find average household income per zip code? Need to find household income per zip code data.",
    2) | Out-Null

# 2. Insert two new bulleted questions right after that paragraph, reusing its
#    list formatting (ListParagraph / numId 4) via InsertParagraphAfter.
$incomeParagraph = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "What income bracket*zip code data.*") {
        $incomeParagraph = $d.Paragraphs.Item($i)
        break
    }
}

$incomeParagraph.Range.InsertParagraphAfter()
$installersParagraph = $incomeParagraph.Next()
$installersParagraph.Range.Text = "Top installers? How much did they install? "

$installersParagraph.Range.InsertParagraphAfter()
$laCountyParagraph = $installersParagraph.Next()
$laCountyParagraph.Range.Text = "Only look at LA county?"

# 3. Merge the two runs split by the old _GoBack bookmark back into a single
#    run ("Why peo" + bookmark + "ple are using solar energy?" -> one run).
$d.Content.Find.Execute(
    "Why people are using solar energy?",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Why people are using solar energy?",
    2) | Out-Null

# 4. Re-anchor the _GoBack bookmark on the now-empty paragraph that follows
#    the two newly-inserted questions (where the old blank paragraph was).
$goBackParagraph = $laCountyParagraph.Next()
$d.Bookmarks.Add("_GoBack", $goBackParagraph.Range) | Out-Null
